$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Harry Potter and the Sorcerers Stone
$ws.Range("A2").Value = "Harry Potter and the Sorcerers Stone"
$ws.Range("B2").Value = "J.K. Rowling"
$ws.Range("C2").Value = Get-Date -Year 2021 -Month 1 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("D2").Value = Get-Date -Year 2021 -Month 1 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("E2").Value = "fiction;wizards;adventure;harry potter"
$ws.Range("F2").Value = "Audio"
$ws.Range("G2").Value = "8 Hours 37 Mins"
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = $true

# Row 3 - Harry Potter and the Chamber of Secrets
$ws.Range("A3").Value = "Harry Potter and the Chamber of Secrets"
$ws.Range("B3").Value = "J.K. Rowling"
$ws.Range("C3").Value = Get-Date -Year 2021 -Month 1 -Day 2 -Hour 0 -Minute 0 -Second 0
$ws.Range("E3").Value = "fiction;wizards;adventure;harry potter"
$ws.Range("F3").Value = "Audio"
$ws.Range("G3").Value = "9 Hours 28 Mins"
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = $true

$ws.Range("E4").Select()
